$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.277.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.780.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +21.69%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.13"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.774.18"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +21.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  +6.04%  "

$ws.Range("E10").Value = "  +8.72%  "

$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.504"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.68"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +12.27%  "

$ws.Range("E14").Value = "  +6.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.393.67"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +21.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.763.05"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +21.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.326.59"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +6.26%  "

$ws.Range("E18").Value = "  +1.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.87%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "521.19"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.96"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +21.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.752"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +9.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +12.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.03"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +8.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.23"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +10.34%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  +10.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.91"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.43"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +15.20%  "

$ws.Range("E33").Value = "  +20.01%  "

$ws.Range("E34").Value = "  +4.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +11.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +10.86%  "

$ws.Range("E38").Value = "  +10.80%  "

$ws.Range("E39").Value = "  +10.32%  "

$ws.Range("E40").Value = "  +9.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "442.09"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +19.09%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.177.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +13.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.95"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.71%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.85"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.18%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0366"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +10.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.60"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.70%  "

$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("E51").Value = "  +8.84%  "
